# Improvements in cast lecture (#21)
#
# 1. Slide 1 ("Type Conversion" title slide): title text gains a trailing "s".
# 2. Slide 3 (Implicit/Explicit conversion diagram): the two callout ovals and
#    the connector pointing at the left oval are swapped/repositioned so the
#    diagram reads "Explicit conversion" on the left (moved slightly right,
#    with a shorter connector) and "Implicit conversion" on the right.
#
# (The deck's masters/layouts also carry a cached "date last saved" field
# that PowerPoint silently refreshes on every save; there is no COM surface
# here that updates that cached field text without collapsing the field into
# plain text, so it is intentionally left alone rather than corrupted.)

function EmuToPoints($emu) {
    # PowerPoint's Left/Top/Width/Height COM properties marshal through a
    # single-precision float expressed in points (1 pt = 12700 EMU); a tiny
    # epsilon keeps values that sit right on a float32 rounding boundary from
    # truncating down to the EMU value just below the intended target.
    return ($emu / 12700.0) + 0.00004
}

$p = $ppt.ActivePresentation

# --- 1. Title slide -------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$titleSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Type Conversions"

# --- 2. Implicit/Explicit conversion diagram slide -------------------------
$diagramSlide = $p.Slides.Item(3)

$leftOval = $diagramSlide.Shapes.Item(2)          # "Oval 7" (id 8)
$leftOval.Left = EmuToPoints 1132216
$leftOval.TextFrame.TextRange.Text = "Explicit conversion"

$connector = $diagramSlide.Shapes.Item(3)         # "Straight Arrow Connector 9" (id 10)
$connector.Left = EmuToPoints 3496236
$connector.Width = EmuToPoints 2599765

$rightOval = $diagramSlide.Shapes.Item(4)         # "Oval 12" (id 13)
$rightOval.TextFrame.TextRange.Text = "Implicit conversion"
